$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start the new row as a copy of the previous data row so that the
# text-typed columns (Date/Weekday/Week) keep their original "typed as
# text" cell type/style instead of being re-interpreted by Excel's
# smart input parsing (which would turn "2024-01-06" into a date serial
# and "00" into the number 0).
$ws.Range("A25:T25").Copy()
$ws.Range("A26:T26").PasteSpecial()

# Now fill in the actual values that changed for this update.
$ws.Cells.Item(26, 2).Value = "20:33:56"

$ws.Cells.Item(26, 5).Value  = 140609
$ws.Cells.Item(26, 6).Value  = 142935
$ws.Cells.Item(26, 7).Value  = 172356
$ws.Cells.Item(26, 8).Value  = 147301
$ws.Cells.Item(26, 9).Value  = -1
$ws.Cells.Item(26, 10).Value = 118359
$ws.Cells.Item(26, 11).Value = 224608
$ws.Cells.Item(26, 12).Value = 249323
$ws.Cells.Item(26, 13).Value = 185077
$ws.Cells.Item(26, 14).Value = 110364
$ws.Cells.Item(26, 15).Value = 40623
$ws.Cells.Item(26, 16).Value = 30810
$ws.Cells.Item(26, 17).Value = 72507
$ws.Cells.Item(26, 18).Value = -1
$ws.Cells.Item(26, 19).Value = 42210
$ws.Cells.Item(26, 20).Value = -1
